$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.861.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.735.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4987"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +7.59%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.18"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07237"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.059"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.15"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.934"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.737.30"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.814"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "86.32"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001033"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06423"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.44"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.731"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "26.929.82"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.25"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.71%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.07"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.933.50"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.207"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.69"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09501"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.581"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.355"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02191"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05854"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.02"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.425"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.15%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.1999"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.763"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6035"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.107"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.622"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.43%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.594"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.60%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5637"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.845"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.74%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.80%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06650"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.08%  "

Write-Host "Update complete"